$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 - this shifts the existing rows 41..119
# down to 42..120 (matching the diff, which shows every row from 41 to 119
# taking on the values previously held by the row below it, and a brand
# new row 120 appearing with the old row 119's data).
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with its data.
$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(41, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(41, 3).Value = "Metropolitana"
$ws.Cells.Item(41, 4).Value = 45014
$ws.Cells.Item(41, 5).Value = 13
$ws.Cells.Item(41, 6).Value = 100114007
$ws.Cells.Item(41, 7).Value = "Jengibre"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 500
$ws.Cells.Item(41, 11).Value = 15000
$ws.Cells.Item(41, 12).Value = 16000
$ws.Cells.Item(41, 13).Value = 15540
$ws.Cells.Item(41, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(41, 15).Value = "Perú"
$ws.Cells.Item(41, 16).Value = 1195
$ws.Cells.Item(41, 17).Value = 13
$ws.Cells.Item(41, 18).Value = "Hortaliza"
